# Corrected Calibration and Ingest Sheets for Coastal Gliders
# - FLORT calibration coefficients on the Asset_Cal_Info sheet:
#     CC_scattering_angle   (row 7, col F): 117   -> 124
#     CC_angular_resolution (row 9, col F): 1.08  -> 1.076
# - Leave the "Asset_Cal_Info" tab as the active/selected sheet (it was the
#   last sheet worked on), with cell C34 selected.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Asset_Cal_Info")

# Make sure we're editing on the right sheet.
$ws.Activate()

# CC_scattering_angle -> 124
$ws.Range("F7").Value = 124

# CC_angular_resolution -> 1.076
$ws.Range("F9").Value = 1.076

# Reflect the author's final cursor position/active tab on save.
$ws.Range("C34").Select()
